$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "studies" - headers + row2 fully rewritten (study metadata columns
# replaced with the new label/description/access_level/contributors/
# reference/reference_year set).
# ---------------------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Cells.Clear()

$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Range("B1").Value = "study_label"
$wsStudies.Range("C1").Value = "description"
$wsStudies.Range("D1").Value = "access_level"
$wsStudies.Range("E1").Value = "contributors"
$wsStudies.Range("F1").Value = "reference"
$wsStudies.Range("G1").Value = "reference_year"

$wsStudies.Range("A2").Value = "foo"
$wsStudies.Range("D2").Value = "public"
$wsStudies.Range("F2").Value = "https://doi.org/10.1093%2Fgenetics%2F16.2.97"
$wsStudies.Range("F2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet "surveys" - columns restructured: study_key -> study_id,
# spatial_notes dropped in favour of location_method/location_notes,
# collection_* columns shift right, and a new time_method column is
# inserted ahead of time_notes.
# ---------------------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Cells.Clear()

$wsSurveys.Range("A1").Value = "study_id"
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("C1").Value = "country_name"
$wsSurveys.Range("D1").Value = "site_name"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("G1").Value = "location_method"
$wsSurveys.Range("H1").Value = "location_notes"
$wsSurveys.Range("H1").Style = "Normal"

$wsSurveys.Range("I1").NumberFormat = "@"
$wsSurveys.Range("I1").Value = "collection_start"
$wsSurveys.Range("J1").NumberFormat = "@"
$wsSurveys.Range("J1").Value = "collection_end"
$wsSurveys.Range("K1").NumberFormat = "@"
$wsSurveys.Range("K1").Value = "collection_day"
$wsSurveys.Range("L1").NumberFormat = "@"
$wsSurveys.Range("L1").Value = "time_method"
$wsSurveys.Range("M1").Value = "time_notes"

$wsSurveys.Range("A2").Value = "foo"
$wsSurveys.Range("B2").Value = "S01"
$wsSurveys.Range("E2").Value = 0
$wsSurveys.Range("F2").Value = 0
$wsSurveys.Range("H2").Value = "example data"
$wsSurveys.Range("H2").Style = "Normal"
$wsSurveys.Range("K2").NumberFormat = "@"
$wsSurveys.Range("K2").Value = "2020-01-01"
$wsSurveys.Range("L2").NumberFormat = "@"
$wsSurveys.Range("M2").Value = "example data"

# ---------------------------------------------------------------------------
# Sheet "counts" - study_key/survey_key headers renamed to study_id/
# survey_id; data rows keep the same variant strings/counts but the study
# key value changes from "study01" to "foo".
# ---------------------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Cells.Clear()

$wsCounts.Range("A1").Value = "study_id"
$wsCounts.Range("A1").Font.Color = 0
$wsCounts.Range("B1").Value = "survey_id"
$wsCounts.Range("C1").Value = "variant_string"
$wsCounts.Range("D1").Value = "variant_num"
$wsCounts.Range("E1").Value = "total_num"

$wsCounts.Range("A2").Value = "foo"
$wsCounts.Range("B2").Value = "S01"
$wsCounts.Range("C2").Value = "crt:1_2_3:AAA;mdr1:1_2_3:AAA"
$wsCounts.Range("D2").Value = 5
$wsCounts.Range("E2").Value = 10

$wsCounts.Range("A3").Value = "foo"
$wsCounts.Range("B3").Value = "S01"
$wsCounts.Range("C3").Value = "crt:1_2_3:AAA;mdr1:1_2_3:AAC"
$wsCounts.Range("D3").Value = 6
$wsCounts.Range("E3").Value = 10

# ---------------------------------------------------------------------------
# Selection / active-tab state. "studies" ends up both the active sheet and
# the one left selected, so its Activate()+Select() run last.
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A3").Select()

$wsCounts.Range("A4").Select()

$wsSurveys.Range("A1:M2").Select()

$wsStudies.Activate()
$wsStudies.Range("D3").Select()
